# Adds a new "REFERENCES" entry below the existing
# "... - pull requests locally" line:
#   https://www.jenkins.io/blog/2016/07/01/html-publisher-plugin/ - publish html
# (as a hyperlink run followed by a plain-text run), matching the commit
# "Modified Jenkinsfile to have reports ran".

$d = $word.ActiveDocument

# Locate the paragraph that ends with " - pull requests locally" - it is the
# anchor after which the new reference paragraph must be inserted.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pull requests locally*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find the 'pull requests locally' reference paragraph"
}

# Insert a brand-new, empty paragraph right after the anchor paragraph.
$tail = $anchor.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# Re-acquire that freshly created paragraph and fill it with plain text first
# (URL + trailing description) so no stray empty run is left behind, then
# turn the URL portion into a hyperlink - mirroring how the existing
# reference paragraphs above it are structured.
$tail.Collapse(0)
$newPara = $tail.Paragraphs(1)

$url = "https://www.jenkins.io/blog/2016/07/01/html-publisher-plugin/"
$suffix = " - publish html"

$newPara.Range.InsertAfter($url + $suffix)

$urlRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $url.Length)
$d.Hyperlinks.Add($urlRange, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null

Write-Output "Inserted reference: $($newPara.Range.Text)"
